# Weekly update: a new Cereza "Lapins / Primera" price record (2023-01-05)
# is inserted above the existing row 144, pushing rows 144-154 down to
# 145-155 (dimension grows from A1:T154 to A1:T155).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 144; Excel copies row 144's
# formatting (incl. the date style on column D) down onto the new row.
$ws.Rows.Item(144).Insert()

$ws.Cells.Item(144, 1).Value2  = 7
$ws.Cells.Item(144, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(144, 3).Value2  = "Ñuble"
$ws.Cells.Item(144, 4).Value2  = 44931
$ws.Cells.Item(144, 5).Value2  = 16
$ws.Cells.Item(144, 6).Value2  = "Fruta"
$ws.Cells.Item(144, 7).Value2  = 100103
$ws.Cells.Item(144, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(144, 9).Value2  = 100103001
$ws.Cells.Item(144, 10).Value2 = "Cereza"
$ws.Cells.Item(144, 11).Value2 = "Lapins"
$ws.Cells.Item(144, 12).Value2 = "Primera"
$ws.Cells.Item(144, 13).Value2 = 160
$ws.Cells.Item(144, 14).Value2 = 4000
$ws.Cells.Item(144, 15).Value2 = 4500
$ws.Cells.Item(144, 16).Value2 = 4250
$ws.Cells.Item(144, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(144, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(144, 19).Value2 = 425
$ws.Cells.Item(144, 20).Value2 = 10
